$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.199924111366272
$ws.Range("B1").Value = 2.606882333755493
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.179608106613159
$ws.Range("E1").Value = 1.172869324684143
